$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$textFormat = "@"

# Row 27
$ws.Cells.Item(27, 1).Value = 112501817
$ws.Cells.Item(27, 2).Value = 56446
$ws.Cells.Item(27, 3).NumberFormat = $textFormat
$ws.Cells.Item(27, 3).Value = "Ovaliderad"
$ws.Cells.Item(27, 4).NumberFormat = $textFormat
$ws.Cells.Item(27, 4).Value = "NT"
$ws.Cells.Item(27, 5).Value = 100049
$ws.Cells.Item(27, 6).NumberFormat = $textFormat
$ws.Cells.Item(27, 6).Value = "Spillkråka"
$ws.Cells.Item(27, 7).NumberFormat = $textFormat
$ws.Cells.Item(27, 7).Value = "Dryocopus martius"
$ws.Cells.Item(27, 8).NumberFormat = $textFormat
$ws.Cells.Item(27, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(27, 9).NumberFormat = $textFormat
$ws.Cells.Item(27, 9).Value = "1"
$ws.Cells.Item(27, 13).NumberFormat = $textFormat
$ws.Cells.Item(27, 13).Value = "födosökande"
$ws.Cells.Item(27, 16).NumberFormat = $textFormat
$ws.Cells.Item(27, 16).Value = "Månses hål, Vstm"
$ws.Cells.Item(27, 17).Value = 558037
$ws.Cells.Item(27, 18).Value = 6628280
$ws.Cells.Item(27, 19).Value = 100
$ws.Cells.Item(27, 20).NumberFormat = $textFormat
$ws.Cells.Item(27, 20).Value = "Västmanland"
$ws.Cells.Item(27, 21).NumberFormat = $textFormat
$ws.Cells.Item(27, 21).Value = "Surahammar"
$ws.Cells.Item(27, 22).NumberFormat = $textFormat
$ws.Cells.Item(27, 22).Value = "Västmanland"
$ws.Cells.Item(27, 23).NumberFormat = $textFormat
$ws.Cells.Item(27, 23).Value = "Ramnäs"
$ws.Cells.Item(27, 25).NumberFormat = $textFormat
$ws.Cells.Item(27, 25).Value = "2023-10-03"
$ws.Cells.Item(27, 26).NumberFormat = $textFormat
$ws.Cells.Item(27, 26).Value = "10:00"
$ws.Cells.Item(27, 27).NumberFormat = $textFormat
$ws.Cells.Item(27, 27).Value = "2023-10-03"
$ws.Cells.Item(27, 28).NumberFormat = $textFormat
$ws.Cells.Item(27, 28).Value = "12:00"
$ws.Cells.Item(27, 30).Value = $false
$ws.Cells.Item(27, 31).Value = $false
$ws.Cells.Item(27, 33).Value = $false
$ws.Cells.Item(27, 49).NumberFormat = $textFormat
$ws.Cells.Item(27, 49).Value = "Tom Sävström"
$ws.Cells.Item(27, 50).NumberFormat = $textFormat
$ws.Cells.Item(27, 50).Value = "Tom Sävström"

# Row 28
$ws.Cells.Item(28, 1).Value = 112501787
$ws.Cells.Item(28, 2).Value = 56575
$ws.Cells.Item(28, 3).NumberFormat = $textFormat
$ws.Cells.Item(28, 3).Value = "Ovaliderad"
$ws.Cells.Item(28, 4).NumberFormat = $textFormat
$ws.Cells.Item(28, 4).Value = "NT"
$ws.Cells.Item(28, 5).Value = 103021
$ws.Cells.Item(28, 6).NumberFormat = $textFormat
$ws.Cells.Item(28, 6).Value = "Talltita"
$ws.Cells.Item(28, 7).NumberFormat = $textFormat
$ws.Cells.Item(28, 7).Value = "Poecile montanus"
$ws.Cells.Item(28, 8).NumberFormat = $textFormat
$ws.Cells.Item(28, 8).Value = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(28, 13).NumberFormat = $textFormat
$ws.Cells.Item(28, 13).Value = "förbiflygande"
$ws.Cells.Item(28, 14).NumberFormat = $textFormat
$ws.Cells.Item(28, 14).Value = "observerad"
$ws.Cells.Item(28, 16).NumberFormat = $textFormat
$ws.Cells.Item(28, 16).Value = "Månses hål, Vstm"
$ws.Cells.Item(28, 17).Value = 558037
$ws.Cells.Item(28, 18).Value = 6628280
$ws.Cells.Item(28, 19).Value = 100
$ws.Cells.Item(28, 20).NumberFormat = $textFormat
$ws.Cells.Item(28, 20).Value = "Västmanland"
$ws.Cells.Item(28, 21).NumberFormat = $textFormat
$ws.Cells.Item(28, 21).Value = "Surahammar"
$ws.Cells.Item(28, 22).NumberFormat = $textFormat
$ws.Cells.Item(28, 22).Value = "Västmanland"
$ws.Cells.Item(28, 23).NumberFormat = $textFormat
$ws.Cells.Item(28, 23).Value = "Ramnäs"
$ws.Cells.Item(28, 25).NumberFormat = $textFormat
$ws.Cells.Item(28, 25).Value = "2023-10-03"
$ws.Cells.Item(28, 26).NumberFormat = $textFormat
$ws.Cells.Item(28, 26).Value = "10:00"
$ws.Cells.Item(28, 27).NumberFormat = $textFormat
$ws.Cells.Item(28, 27).Value = "2023-10-03"
$ws.Cells.Item(28, 28).NumberFormat = $textFormat
$ws.Cells.Item(28, 28).Value = "12:00"
$ws.Cells.Item(28, 29).NumberFormat = $textFormat
$ws.Cells.Item(28, 29).Value = "meståg"
$ws.Cells.Item(28, 30).Value = $false
$ws.Cells.Item(28, 31).Value = $false
$ws.Cells.Item(28, 33).Value = $false
$ws.Cells.Item(28, 49).NumberFormat = $textFormat
$ws.Cells.Item(28, 49).Value = "Tom Sävström"
$ws.Cells.Item(28, 50).NumberFormat = $textFormat
$ws.Cells.Item(28, 50).Value = "Tom Sävström"

# Row 29
$ws.Cells.Item(29, 1).Value = 112501637
$ws.Cells.Item(29, 2).Value = 90814
$ws.Cells.Item(29, 3).NumberFormat = $textFormat
$ws.Cells.Item(29, 3).Value = "Ovaliderad"
$ws.Cells.Item(29, 4).NumberFormat = $textFormat
$ws.Cells.Item(29, 4).Value = "LC"
$ws.Cells.Item(29, 5).Value = 4364
$ws.Cells.Item(29, 6).NumberFormat = $textFormat
$ws.Cells.Item(29, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(29, 7).NumberFormat = $textFormat
$ws.Cells.Item(29, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(29, 8).NumberFormat = $textFormat
$ws.Cells.Item(29, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(29, 16).NumberFormat = $textFormat
$ws.Cells.Item(29, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(29, 17).Value = 557914
$ws.Cells.Item(29, 18).Value = 6628494
$ws.Cells.Item(29, 19).Value = 10
$ws.Cells.Item(29, 20).NumberFormat = $textFormat
$ws.Cells.Item(29, 20).Value = "Västmanland"
$ws.Cells.Item(29, 21).NumberFormat = $textFormat
$ws.Cells.Item(29, 21).Value = "Surahammar"
$ws.Cells.Item(29, 22).NumberFormat = $textFormat
$ws.Cells.Item(29, 22).Value = "Västmanland"
$ws.Cells.Item(29, 23).NumberFormat = $textFormat
$ws.Cells.Item(29, 23).Value = "Ramnäs"
$ws.Cells.Item(29, 25).NumberFormat = $textFormat
$ws.Cells.Item(29, 25).Value = "2023-10-03"
$ws.Cells.Item(29, 27).NumberFormat = $textFormat
$ws.Cells.Item(29, 27).Value = "2023-10-03"
$ws.Cells.Item(29, 30).Value = $false
$ws.Cells.Item(29, 31).Value = $false
$ws.Cells.Item(29, 33).Value = $false
$ws.Cells.Item(29, 35).NumberFormat = $textFormat
$ws.Cells.Item(29, 35).Value = "Barrblandskog"
$ws.Cells.Item(29, 49).NumberFormat = $textFormat
$ws.Cells.Item(29, 49).Value = "Tom Sävström"
$ws.Cells.Item(29, 50).NumberFormat = $textFormat
$ws.Cells.Item(29, 50).Value = "Tom Sävström"

# Row 30
$ws.Cells.Item(30, 1).Value = 112501403
$ws.Cells.Item(30, 2).Value = 96735
$ws.Cells.Item(30, 3).NumberFormat = $textFormat
$ws.Cells.Item(30, 3).Value = "Ovaliderad"
$ws.Cells.Item(30, 4).NumberFormat = $textFormat
$ws.Cells.Item(30, 4).Value = "VU"
$ws.Cells.Item(30, 5).Value = 220787
$ws.Cells.Item(30, 6).NumberFormat = $textFormat
$ws.Cells.Item(30, 6).Value = "Knärot"
$ws.Cells.Item(30, 7).NumberFormat = $textFormat
$ws.Cells.Item(30, 7).Value = "Goodyera repens"
$ws.Cells.Item(30, 8).NumberFormat = $textFormat
$ws.Cells.Item(30, 8).Value = "(L.) R. Br."
$ws.Cells.Item(30, 9).NumberFormat = $textFormat
$ws.Cells.Item(30, 9).Value = "14"
$ws.Cells.Item(30, 10).NumberFormat = $textFormat
$ws.Cells.Item(30, 10).Value = "plantor/tuvor"
$ws.Cells.Item(30, 11).NumberFormat = $textFormat
$ws.Cells.Item(30, 11).Value = "fullt utvecklade blad"
$ws.Cells.Item(30, 16).NumberFormat = $textFormat
$ws.Cells.Item(30, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(30, 17).Value = 557958
$ws.Cells.Item(30, 18).Value = 6628535
$ws.Cells.Item(30, 19).Value = 10
$ws.Cells.Item(30, 20).NumberFormat = $textFormat
$ws.Cells.Item(30, 20).Value = "Västmanland"
$ws.Cells.Item(30, 21).NumberFormat = $textFormat
$ws.Cells.Item(30, 21).Value = "Surahammar"
$ws.Cells.Item(30, 22).NumberFormat = $textFormat
$ws.Cells.Item(30, 22).Value = "Västmanland"
$ws.Cells.Item(30, 23).NumberFormat = $textFormat
$ws.Cells.Item(30, 23).Value = "Ramnäs"
$ws.Cells.Item(30, 25).NumberFormat = $textFormat
$ws.Cells.Item(30, 25).Value = "2023-10-03"
$ws.Cells.Item(30, 27).NumberFormat = $textFormat
$ws.Cells.Item(30, 27).Value = "2023-10-03"
$ws.Cells.Item(30, 30).Value = $false
$ws.Cells.Item(30, 31).Value = $false
$ws.Cells.Item(30, 33).Value = $false
$ws.Cells.Item(30, 35).NumberFormat = $textFormat
$ws.Cells.Item(30, 35).Value = "Barrblandskog"
$ws.Cells.Item(30, 49).NumberFormat = $textFormat
$ws.Cells.Item(30, 49).Value = "Tom Sävström"
$ws.Cells.Item(30, 50).NumberFormat = $textFormat
$ws.Cells.Item(30, 50).Value = "Tom Sävström"

# Row 31
$ws.Cells.Item(31, 1).Value = 112501672
$ws.Cells.Item(31, 2).Value = 89517
$ws.Cells.Item(31, 3).NumberFormat = $textFormat
$ws.Cells.Item(31, 3).Value = "Ovaliderad"
$ws.Cells.Item(31, 4).NumberFormat = $textFormat
$ws.Cells.Item(31, 4).Value = "LC"
$ws.Cells.Item(31, 5).Value = 5447
$ws.Cells.Item(31, 6).NumberFormat = $textFormat
$ws.Cells.Item(31, 6).Value = "Vedticka"
$ws.Cells.Item(31, 7).NumberFormat = $textFormat
$ws.Cells.Item(31, 7).Value = "Fuscoporia viticola"
$ws.Cells.Item(31, 8).NumberFormat = $textFormat
$ws.Cells.Item(31, 8).Value = "(Schwein.) Murrill"
$ws.Cells.Item(31, 16).NumberFormat = $textFormat
$ws.Cells.Item(31, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(31, 17).Value = 557941
$ws.Cells.Item(31, 18).Value = 6628426
$ws.Cells.Item(31, 19).Value = 10
$ws.Cells.Item(31, 20).NumberFormat = $textFormat
$ws.Cells.Item(31, 20).Value = "Västmanland"
$ws.Cells.Item(31, 21).NumberFormat = $textFormat
$ws.Cells.Item(31, 21).Value = "Surahammar"
$ws.Cells.Item(31, 22).NumberFormat = $textFormat
$ws.Cells.Item(31, 22).Value = "Västmanland"
$ws.Cells.Item(31, 23).NumberFormat = $textFormat
$ws.Cells.Item(31, 23).Value = "Ramnäs"
$ws.Cells.Item(31, 25).NumberFormat = $textFormat
$ws.Cells.Item(31, 25).Value = "2023-10-03"
$ws.Cells.Item(31, 27).NumberFormat = $textFormat
$ws.Cells.Item(31, 27).Value = "2023-10-03"
$ws.Cells.Item(31, 30).Value = $false
$ws.Cells.Item(31, 31).Value = $false
$ws.Cells.Item(31, 33).Value = $false
$ws.Cells.Item(31, 35).NumberFormat = $textFormat
$ws.Cells.Item(31, 35).Value = "Barrblandskog"
$ws.Cells.Item(31, 41).NumberFormat = $textFormat
$ws.Cells.Item(31, 41).Value = "Gran"
$ws.Cells.Item(31, 49).NumberFormat = $textFormat
$ws.Cells.Item(31, 49).Value = "Tom Sävström"
$ws.Cells.Item(31, 50).NumberFormat = $textFormat
$ws.Cells.Item(31, 50).Value = "Tom Sävström"

# Row 32
$ws.Cells.Item(32, 1).Value = 112501295
$ws.Cells.Item(32, 2).Value = 8377
$ws.Cells.Item(32, 3).NumberFormat = $textFormat
$ws.Cells.Item(32, 3).Value = "Ovaliderad"
$ws.Cells.Item(32, 4).NumberFormat = $textFormat
$ws.Cells.Item(32, 4).Value = "LC"
$ws.Cells.Item(32, 5).Value = 106545
$ws.Cells.Item(32, 6).NumberFormat = $textFormat
$ws.Cells.Item(32, 6).Value = "Mindre märgborre"
$ws.Cells.Item(32, 7).NumberFormat = $textFormat
$ws.Cells.Item(32, 7).Value = "Tomicus minor"
$ws.Cells.Item(32, 8).NumberFormat = $textFormat
$ws.Cells.Item(32, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(32, 13).NumberFormat = $textFormat
$ws.Cells.Item(32, 13).Value = "äldre gnagspår"
$ws.Cells.Item(32, 16).NumberFormat = $textFormat
$ws.Cells.Item(32, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(32, 17).Value = 558049
$ws.Cells.Item(32, 18).Value = 6628477
$ws.Cells.Item(32, 19).Value = 10
$ws.Cells.Item(32, 20).NumberFormat = $textFormat
$ws.Cells.Item(32, 20).Value = "Västmanland"
$ws.Cells.Item(32, 21).NumberFormat = $textFormat
$ws.Cells.Item(32, 21).Value = "Surahammar"
$ws.Cells.Item(32, 22).NumberFormat = $textFormat
$ws.Cells.Item(32, 22).Value = "Västmanland"
$ws.Cells.Item(32, 23).NumberFormat = $textFormat
$ws.Cells.Item(32, 23).Value = "Ramnäs"
$ws.Cells.Item(32, 25).NumberFormat = $textFormat
$ws.Cells.Item(32, 25).Value = "2023-10-03"
$ws.Cells.Item(32, 27).NumberFormat = $textFormat
$ws.Cells.Item(32, 27).Value = "2023-10-03"
$ws.Cells.Item(32, 30).Value = $false
$ws.Cells.Item(32, 31).Value = $false
$ws.Cells.Item(32, 33).Value = $false
$ws.Cells.Item(32, 35).NumberFormat = $textFormat
$ws.Cells.Item(32, 35).Value = "Tallskog, inslag av unga granplantor"
$ws.Cells.Item(32, 41).NumberFormat = $textFormat
$ws.Cells.Item(32, 41).Value = "Tall"
$ws.Cells.Item(32, 49).NumberFormat = $textFormat
$ws.Cells.Item(32, 49).Value = "Tom Sävström"
$ws.Cells.Item(32, 50).NumberFormat = $textFormat
$ws.Cells.Item(32, 50).Value = "Tom Sävström"

# Row 33
$ws.Cells.Item(33, 1).Value = 112501432
$ws.Cells.Item(33, 2).Value = 89517
$ws.Cells.Item(33, 3).NumberFormat = $textFormat
$ws.Cells.Item(33, 3).Value = "Ovaliderad"
$ws.Cells.Item(33, 4).NumberFormat = $textFormat
$ws.Cells.Item(33, 4).Value = "LC"
$ws.Cells.Item(33, 5).Value = 5447
$ws.Cells.Item(33, 6).NumberFormat = $textFormat
$ws.Cells.Item(33, 6).Value = "Vedticka"
$ws.Cells.Item(33, 7).NumberFormat = $textFormat
$ws.Cells.Item(33, 7).Value = "Fuscoporia viticola"
$ws.Cells.Item(33, 8).NumberFormat = $textFormat
$ws.Cells.Item(33, 8).Value = "(Schwein.) Murrill"
$ws.Cells.Item(33, 16).NumberFormat = $textFormat
$ws.Cells.Item(33, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(33, 17).Value = 557940
$ws.Cells.Item(33, 18).Value = 6628530
$ws.Cells.Item(33, 19).Value = 10
$ws.Cells.Item(33, 20).NumberFormat = $textFormat
$ws.Cells.Item(33, 20).Value = "Västmanland"
$ws.Cells.Item(33, 21).NumberFormat = $textFormat
$ws.Cells.Item(33, 21).Value = "Surahammar"
$ws.Cells.Item(33, 22).NumberFormat = $textFormat
$ws.Cells.Item(33, 22).Value = "Västmanland"
$ws.Cells.Item(33, 23).NumberFormat = $textFormat
$ws.Cells.Item(33, 23).Value = "Ramnäs"
$ws.Cells.Item(33, 25).NumberFormat = $textFormat
$ws.Cells.Item(33, 25).Value = "2023-10-03"
$ws.Cells.Item(33, 27).NumberFormat = $textFormat
$ws.Cells.Item(33, 27).Value = "2023-10-03"
$ws.Cells.Item(33, 30).Value = $false
$ws.Cells.Item(33, 31).Value = $false
$ws.Cells.Item(33, 33).Value = $false
$ws.Cells.Item(33, 35).NumberFormat = $textFormat
$ws.Cells.Item(33, 35).Value = "Barrblandskog"
$ws.Cells.Item(33, 41).NumberFormat = $textFormat
$ws.Cells.Item(33, 41).Value = "Gran"
$ws.Cells.Item(33, 49).NumberFormat = $textFormat
$ws.Cells.Item(33, 49).Value = "Tom Sävström"
$ws.Cells.Item(33, 50).NumberFormat = $textFormat
$ws.Cells.Item(33, 50).Value = "Tom Sävström"

# Row 34
$ws.Cells.Item(34, 1).Value = 112501340
$ws.Cells.Item(34, 2).Value = 90814
$ws.Cells.Item(34, 3).NumberFormat = $textFormat
$ws.Cells.Item(34, 3).Value = "Ovaliderad"
$ws.Cells.Item(34, 4).NumberFormat = $textFormat
$ws.Cells.Item(34, 4).Value = "LC"
$ws.Cells.Item(34, 5).Value = 4364
$ws.Cells.Item(34, 6).NumberFormat = $textFormat
$ws.Cells.Item(34, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(34, 7).NumberFormat = $textFormat
$ws.Cells.Item(34, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(34, 8).NumberFormat = $textFormat
$ws.Cells.Item(34, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(34, 16).NumberFormat = $textFormat
$ws.Cells.Item(34, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(34, 17).Value = 557994
$ws.Cells.Item(34, 18).Value = 6628516
$ws.Cells.Item(34, 19).Value = 10
$ws.Cells.Item(34, 20).NumberFormat = $textFormat
$ws.Cells.Item(34, 20).Value = "Västmanland"
$ws.Cells.Item(34, 21).NumberFormat = $textFormat
$ws.Cells.Item(34, 21).Value = "Surahammar"
$ws.Cells.Item(34, 22).NumberFormat = $textFormat
$ws.Cells.Item(34, 22).Value = "Västmanland"
$ws.Cells.Item(34, 23).NumberFormat = $textFormat
$ws.Cells.Item(34, 23).Value = "Ramnäs"
$ws.Cells.Item(34, 25).NumberFormat = $textFormat
$ws.Cells.Item(34, 25).Value = "2023-10-03"
$ws.Cells.Item(34, 27).NumberFormat = $textFormat
$ws.Cells.Item(34, 27).Value = "2023-10-03"
$ws.Cells.Item(34, 30).Value = $false
$ws.Cells.Item(34, 31).Value = $false
$ws.Cells.Item(34, 33).Value = $false
$ws.Cells.Item(34, 35).NumberFormat = $textFormat
$ws.Cells.Item(34, 35).Value = "Barrblandskog, mot tallmosse"
$ws.Cells.Item(34, 49).NumberFormat = $textFormat
$ws.Cells.Item(34, 49).Value = "Tom Sävström"
$ws.Cells.Item(34, 50).NumberFormat = $textFormat
$ws.Cells.Item(34, 50).Value = "Tom Sävström"

# Row 35
$ws.Cells.Item(35, 1).Value = 112501658
$ws.Cells.Item(35, 2).Value = 96735
$ws.Cells.Item(35, 3).NumberFormat = $textFormat
$ws.Cells.Item(35, 3).Value = "Ovaliderad"
$ws.Cells.Item(35, 4).NumberFormat = $textFormat
$ws.Cells.Item(35, 4).Value = "VU"
$ws.Cells.Item(35, 5).Value = 220787
$ws.Cells.Item(35, 6).NumberFormat = $textFormat
$ws.Cells.Item(35, 6).Value = "Knärot"
$ws.Cells.Item(35, 7).NumberFormat = $textFormat
$ws.Cells.Item(35, 7).Value = "Goodyera repens"
$ws.Cells.Item(35, 8).NumberFormat = $textFormat
$ws.Cells.Item(35, 8).Value = "(L.) R. Br."
$ws.Cells.Item(35, 9).NumberFormat = $textFormat
$ws.Cells.Item(35, 9).Value = "16"
$ws.Cells.Item(35, 10).NumberFormat = $textFormat
$ws.Cells.Item(35, 10).Value = "plantor/tuvor"
$ws.Cells.Item(35, 11).NumberFormat = $textFormat
$ws.Cells.Item(35, 11).Value = "fullt utvecklade blad"
$ws.Cells.Item(35, 16).NumberFormat = $textFormat
$ws.Cells.Item(35, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(35, 17).Value = 557921
$ws.Cells.Item(35, 18).Value = 6628494
$ws.Cells.Item(35, 19).Value = 10
$ws.Cells.Item(35, 20).NumberFormat = $textFormat
$ws.Cells.Item(35, 20).Value = "Västmanland"
$ws.Cells.Item(35, 21).NumberFormat = $textFormat
$ws.Cells.Item(35, 21).Value = "Surahammar"
$ws.Cells.Item(35, 22).NumberFormat = $textFormat
$ws.Cells.Item(35, 22).Value = "Västmanland"
$ws.Cells.Item(35, 23).NumberFormat = $textFormat
$ws.Cells.Item(35, 23).Value = "Ramnäs"
$ws.Cells.Item(35, 25).NumberFormat = $textFormat
$ws.Cells.Item(35, 25).Value = "2023-10-03"
$ws.Cells.Item(35, 27).NumberFormat = $textFormat
$ws.Cells.Item(35, 27).Value = "2023-10-03"
$ws.Cells.Item(35, 30).Value = $false
$ws.Cells.Item(35, 31).Value = $false
$ws.Cells.Item(35, 33).Value = $false
$ws.Cells.Item(35, 35).NumberFormat = $textFormat
$ws.Cells.Item(35, 35).Value = "Barrblandskog"
$ws.Cells.Item(35, 49).NumberFormat = $textFormat
$ws.Cells.Item(35, 49).Value = "Tom Sävström"
$ws.Cells.Item(35, 50).NumberFormat = $textFormat
$ws.Cells.Item(35, 50).Value = "Tom Sävström"

# Row 36
$ws.Cells.Item(36, 1).Value = 112501410
$ws.Cells.Item(36, 2).Value = 89517
$ws.Cells.Item(36, 3).NumberFormat = $textFormat
$ws.Cells.Item(36, 3).Value = "Ovaliderad"
$ws.Cells.Item(36, 4).NumberFormat = $textFormat
$ws.Cells.Item(36, 4).Value = "LC"
$ws.Cells.Item(36, 5).Value = 5447
$ws.Cells.Item(36, 6).NumberFormat = $textFormat
$ws.Cells.Item(36, 6).Value = "Vedticka"
$ws.Cells.Item(36, 7).NumberFormat = $textFormat
$ws.Cells.Item(36, 7).Value = "Fuscoporia viticola"
$ws.Cells.Item(36, 8).NumberFormat = $textFormat
$ws.Cells.Item(36, 8).Value = "(Schwein.) Murrill"
$ws.Cells.Item(36, 16).NumberFormat = $textFormat
$ws.Cells.Item(36, 16).Value = "Månses hål N, Vstm"
$ws.Cells.Item(36, 17).Value = 557958
$ws.Cells.Item(36, 18).Value = 6628535
$ws.Cells.Item(36, 19).Value = 10
$ws.Cells.Item(36, 20).NumberFormat = $textFormat
$ws.Cells.Item(36, 20).Value = "Västmanland"
$ws.Cells.Item(36, 21).NumberFormat = $textFormat
$ws.Cells.Item(36, 21).Value = "Surahammar"
$ws.Cells.Item(36, 22).NumberFormat = $textFormat
$ws.Cells.Item(36, 22).Value = "Västmanland"
$ws.Cells.Item(36, 23).NumberFormat = $textFormat
$ws.Cells.Item(36, 23).Value = "Ramnäs"
$ws.Cells.Item(36, 25).NumberFormat = $textFormat
$ws.Cells.Item(36, 25).Value = "2023-10-03"
$ws.Cells.Item(36, 27).NumberFormat = $textFormat
$ws.Cells.Item(36, 27).Value = "2023-10-03"
$ws.Cells.Item(36, 30).Value = $false
$ws.Cells.Item(36, 31).Value = $false
$ws.Cells.Item(36, 33).Value = $false
$ws.Cells.Item(36, 35).NumberFormat = $textFormat
$ws.Cells.Item(36, 35).Value = "Barrblandskog"
$ws.Cells.Item(36, 41).NumberFormat = $textFormat
$ws.Cells.Item(36, 41).Value = "Gran"
$ws.Cells.Item(36, 49).NumberFormat = $textFormat
$ws.Cells.Item(36, 49).Value = "Tom Sävström"
$ws.Cells.Item(36, 50).NumberFormat = $textFormat
$ws.Cells.Item(36, 50).Value = "Tom Sävström"

